# Update the summary ("Сводная таблица") pivot table:
#  - rename the "ММП v2" algorithm label to "ММП" everywhere
#  - refresh the recomputed metric values (Nf/Nf_max/Iter/Iter_max/t/t_max/
#    Рэф1-4/sqrt(Ка)) produced by the updated/dynamic algorithm run.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2, 1).Range.Text = "ММП"
$t.Cell(2, 8).Range.Text = "12"
$t.Cell(2, 9).Range.Text = "26"
$t.Cell(2, 10).Range.Text = "10"
$t.Cell(2, 11).Range.Text = "20"
$t.Cell(2, 12).Range.Text = "0.0031"
$t.Cell(2, 13).Range.Text = "0.0113"
$t.Cell(3, 12).Range.Text = "0.0003"
$t.Cell(4, 1).Range.Text = "ММП"
$t.Cell(4, 3).Range.Text = "0.96"
$t.Cell(4, 4).Range.Text = "1.0"
$t.Cell(4, 5).Range.Text = "1.0"
$t.Cell(4, 6).Range.Text = "1.0"
$t.Cell(4, 9).Range.Text = "26"
$t.Cell(4, 10).Range.Text = "10"
$t.Cell(4, 11).Range.Text = "20"
$t.Cell(4, 12).Range.Text = "0.0028"
$t.Cell(4, 13).Range.Text = "0.0060"
$t.Cell(5, 3).Range.Text = "0.68"
$t.Cell(5, 4).Range.Text = "0.88"
$t.Cell(5, 5).Range.Text = "0.88"
$t.Cell(5, 6).Range.Text = "0.88"
$t.Cell(5, 12).Range.Text = "0.0003"
$t.Cell(5, 13).Range.Text = "0.0004"
$t.Cell(6, 1).Range.Text = "ММП"
$t.Cell(6, 3).Range.Text = "0.92"
$t.Cell(6, 4).Range.Text = "0.96"
$t.Cell(6, 5).Range.Text = "0.96"
$t.Cell(6, 6).Range.Text = "0.96"
$t.Cell(6, 8).Range.Text = "12"
$t.Cell(6, 9).Range.Text = "26"
$t.Cell(6, 10).Range.Text = "10"
$t.Cell(6, 11).Range.Text = "20"
$t.Cell(6, 12).Range.Text = "0.0028"
$t.Cell(6, 13).Range.Text = "0.0067"
$t.Cell(7, 3).Range.Text = "0.36"
$t.Cell(7, 4).Range.Text = "0.66"
$t.Cell(7, 5).Range.Text = "0.66"
$t.Cell(7, 6).Range.Text = "0.74"
$t.Cell(7, 12).Range.Text = "0.0003"
$t.Cell(7, 13).Range.Text = "0.0005"
$t.Cell(8, 1).Range.Text = "ММП"
$t.Cell(8, 3).Range.Text = "0.82"
$t.Cell(8, 4).Range.Text = "0.94"
$t.Cell(8, 5).Range.Text = "0.94"
$t.Cell(8, 6).Range.Text = "0.94"
$t.Cell(8, 9).Range.Text = "26"
$t.Cell(8, 10).Range.Text = "10"
$t.Cell(8, 11).Range.Text = "20"
$t.Cell(8, 12).Range.Text = "0.0028"
$t.Cell(8, 13).Range.Text = "0.0061"
$t.Cell(9, 3).Range.Text = "0.26"
$t.Cell(9, 4).Range.Text = "0.46"
$t.Cell(9, 5).Range.Text = "0.5"
$t.Cell(9, 6).Range.Text = "0.56"
$t.Cell(9, 12).Range.Text = "0.0003"
$t.Cell(9, 13).Range.Text = "0.0007"
$t.Cell(10, 1).Range.Text = "ММП"
$t.Cell(10, 3).Range.Text = "0.68"
$t.Cell(10, 4).Range.Text = "0.88"
$t.Cell(10, 5).Range.Text = "0.88"
$t.Cell(10, 6).Range.Text = "0.94"
$t.Cell(10, 9).Range.Text = "26"
$t.Cell(10, 10).Range.Text = "10"
$t.Cell(10, 11).Range.Text = "20"
$t.Cell(10, 12).Range.Text = "0.0027"
$t.Cell(10, 13).Range.Text = "0.0057"
$t.Cell(11, 3).Range.Text = "0.12"
$t.Cell(11, 4).Range.Text = "0.24"
$t.Cell(11, 5).Range.Text = "0.24"
$t.Cell(11, 6).Range.Text = "0.3"
$t.Cell(11, 7).Range.Text = "0.60"
$t.Cell(11, 12).Range.Text = "0.0003"
$t.Cell(11, 13).Range.Text = "0.0014"
$t.Cell(12, 1).Range.Text = "ММП"
$t.Cell(12, 3).Range.Text = "0.36"
$t.Cell(12, 4).Range.Text = "0.7"
$t.Cell(12, 5).Range.Text = "0.72"
$t.Cell(12, 6).Range.Text = "0.76"
$t.Cell(12, 8).Range.Text = "12"
$t.Cell(12, 9).Range.Text = "26"
$t.Cell(12, 10).Range.Text = "10"
$t.Cell(12, 11).Range.Text = "20"
$t.Cell(12, 12).Range.Text = "0.0028"
$t.Cell(12, 13).Range.Text = "0.0062"
$t.Cell(13, 3).Range.Text = "0.0"
$t.Cell(13, 4).Range.Text = "0.0"
$t.Cell(13, 5).Range.Text = "0.0"
$t.Cell(13, 6).Range.Text = "0.0"
$t.Cell(13, 7).Range.Text = "1.77"
$t.Cell(13, 12).Range.Text = "0.0003"
$t.Cell(13, 13).Range.Text = "0.0012"
